$d = $word.ActiveDocument
$d.Content.Find.Execute("Conatndo", $true, $false, $false, $false, $false, $true, 1, $false, "Contando", 2)
$d.Content.Find.Execute("Comntando", $true, $false, $false, $false, $false, $true, 1, $false, "Contando", 2)
